# Applies the "Updated cryptos list" data refresh to sheet1 of the workbook.
# For each changed cell we assign the new literal text value. Cells in column D
# that look like plain decimal numbers (e.g. "599.02") are written with a leading
# apostrophe so Excel stores them as text (matching the source file's inlineStr
# cells) instead of silently converting them to numeric values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.923.69'
$ws.Range("E2").Value = '  -2.43%  '

$ws.Range("D3").Value = '3.186.18'
$ws.Range("E3").Value = '  -1.39%  '

$ws.Range("E4").Value = '  +0.08%  '

$ws.Range("D5").Value = '''599.02'

$ws.Range("D6").Value = '''153.04'
$ws.Range("E6").Value = '  -3.26%  '

$ws.Range("E7").Value = '  +0.05%  '

$ws.Range("D8").Value = '3.183.56'
$ws.Range("E8").Value = '  -1.46%  '

$ws.Range("D9").Value = '''0.527'
$ws.Range("E9").Value = '  -3.87%  '

$ws.Range("E10").Value = '  -4.58%  '

$ws.Range("D11").Value = '''5.59'
$ws.Range("E11").Value = '  -1.71%  '

$ws.Range("D12").Value = '''0.479'
$ws.Range("E12").Value = '  -5.79%  '

$ws.Range("D13").Value = '''0.0000259'
$ws.Range("E13").Value = '  -5.55%  '

$ws.Range("D14").Value = '''37.13'
$ws.Range("E14").Value = '  -4.96%  '

$ws.Range("D15").Value = '3.719.90'
$ws.Range("E15").Value = '  -1.12%  '

$ws.Range("D16").Value = '64.945.13'
$ws.Range("E16").Value = '  -2.52%  '

$ws.Range("D17").Value = '3.194.60'
$ws.Range("E17").Value = '  -1.15%  '

$ws.Range("E18").Value = '  +0.35%  '

$ws.Range("D19").Value = '''7.06'
$ws.Range("E19").Value = '  -5.11%  '

$ws.Range("D20").Value = '''484.19'
$ws.Range("E20").Value = '  -5.29%  '

$ws.Range("D21").Value = '''14.83'
$ws.Range("E21").Value = '  -3.07%  '

$ws.Range("D22").Value = '''0.719'
$ws.Range("E22").Value = '  -2.36%  '

$ws.Range("D23").Value = '''7.77'
$ws.Range("E23").Value = '  -3.89%  '

$ws.Range("D24").Value = '''13.93'
$ws.Range("E24").Value = '  -5.69%  '

$ws.Range("D25").Value = '''85.42'
$ws.Range("E25").Value = '  +0.70%  '

$ws.Range("E26").Value = '  -0.18%  '

$ws.Range("E27").Value = '  -2.02%  '

$ws.Range("D28").Value = '''8.70'
$ws.Range("E28").Value = '  -5.12%  '

# Rows 29 and 30: two coins swap rank position (plus updated price/volume)
$ws.Range("B29").Value = 'Hedera'
$ws.Range("C29").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D29").Value = '''0.128'
$ws.Range("E29").Value = '  +32.22%  '

$ws.Range("B30").Value = 'ImmutableX'
$ws.Range("C30").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D30").Value = '''2.28'
$ws.Range("E30").Value = '  -5.24%  '

$ws.Range("D31").Value = '''6.98'
$ws.Range("E31").Value = '  -1.09%  '

$ws.Range("E32").Value = '  -8.71%  '

# Rows 33 and 34: two coins swap rank position (plus updated price/volume)
$ws.Range("B33").Value = 'FirstDigitalUSD'
$ws.Range("C33").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D33").Value = '''1.00'
$ws.Range("E33").Value = '  -0.12%  '

$ws.Range("B34").Value = 'EthereumClassic'
$ws.Range("C34").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D34").Value = '''26.96'
$ws.Range("E34").Value = '  -4.60%  '

$ws.Range("E35").Value = '  -6.28%  '

$ws.Range("D36").Value = '''6.14'
$ws.Range("E36").Value = '  -5.88%  '

$ws.Range("D37").Value = '''3.30'
$ws.Range("E37").Value = '  +8.37%  '

$ws.Range("D38").Value = '''54.51'
$ws.Range("E38").Value = '  -2.93%  '

$ws.Range("D39").Value = '''476.47'

$ws.Range("D40").Value = '0.0₃0730'
$ws.Range("E40").Value = '  -5.77%  '

$ws.Range("D41").Value = '''0.0405'
$ws.Range("E41").Value = '  -3.89%  '

$ws.Range("E42").Value = '  -2.89%  '

$ws.Range("D43").Value = '''8.54'
$ws.Range("E43").Value = '  -2.87%  '

$ws.Range("D44").Value = '2.917.41'
$ws.Range("E44").Value = '  +1.27%  '

$ws.Range("D45").Value = '''2.44'
$ws.Range("E45").Value = '  -1.93%  '

$ws.Range("D46").Value = '''0.277'
$ws.Range("E46").Value = '  -7.61%  '

$ws.Range("D47").Value = '''27.53'
$ws.Range("E47").Value = '  -3.76%  '

# Rows 48 and 49: two coins swap rank position (plus updated price/volume)
$ws.Range("B48").Value = 'USDe'
$ws.Range("C48").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D48").Value = '''0.998'
$ws.Range("E48").Value = '  -0.05%  '

$ws.Range("B49").Value = 'ThetaToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D49").Value = '''2.36'
$ws.Range("E49").Value = '  -2.06%  '

$ws.Range("D50").Value = '''0.116'
$ws.Range("E50").Value = '  -0.43%  '

$ws.Range("D51").Value = '''121.13'
$ws.Range("E51").Value = '  -1.76%  '
